# Update the GitHub link on the "Getting started" slide so the
# "www.github.com/stevenm1/ez" text becomes a clickable hyperlink.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$fullText = $tr.Text
$target = "www.github.com/stevenm1/ez"
$idx = $fullText.IndexOf($target)

if ($idx -ge 0) {
    # TextRange.Characters(Start, Length) is 1-based.
    $linkRange = $tr.Characters($idx + 1, $target.Length)
    $action = $linkRange.ActionSettings(1)
    $action.Hyperlink.Address = "https://www.github.com/stevenm1/ez"
}
